{"js": "// Trim trailing whitespace from several Q&A runs in the \"preguntas gen\u00e9ricas\"\n// section of the document (per the commit's underlying OOXML diff).\n//\n// Each edit below:\n//   1. Searches for the exact text (including the trailing space that must go).\n//   2. Replaces that matched range's text with the trimmed version.\n//\n// Using Range.search + Range.insertText(\"Replace\") only rewrites the text of\n// the run(s) actually covered by the match, which keeps the edit minimal and\n// matches the target diff (no unrelated runs/paragraphs are touched).\n\nasync function trimTrailingSpace(context, searchText, replacementText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 1. \"...dependiendo el tipo de trabajo realizado se tienen diferentes costos \"\nawait trimTrailingSpace(\n  context,\n  \" dependiendo el tipo de trabajo realizado se tienen diferentes costos \",\n  \" dependiendo el tipo de trabajo realizado se tienen diferentes costos\"\n);\n\n// 2. \"video conferencia \"\nawait trimTrailingSpace(\n  context,\n  \"video conferencia \",\n  \"video conferencia\"\n);\n\n// 3. \"8. \u00bfCu\u00e1l es su prop\u00f3sito? \" \u2014 drop the trailing space-only run after \"prop\u00f3sito?\"\nawait trimTrailingSpace(\n  context,\n  \"prop\u00f3sito? \",\n  \"prop\u00f3sito?\"\n);\n\n// 4. \"Que necesito para realizar la consultor\u00eda? \"\nawait trimTrailingSpace(\n  context,\n  \"Que necesito para realizar la consultor\u00eda? \",\n  \"Que necesito para realizar la consultor\u00eda?\"\n);\n\n// 5. \"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor? \"\nawait trimTrailingSpace(\n  context,\n  \"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor? \",\n  \"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor?\"\n);\n\n// 6. \"agenda otra cita \"\nawait trimTrailingSpace(\n  context,\n  \"agenda otra cita \",\n  \"agenda otra cita\"\n);\n", "ps1": "# Trim trailing whitespace from several Q&A runs in the \"preguntas gen\u00e9ricas\"\n# section of the document (per the commit's underlying OOXML diff).\n#\n# For each target phrase we locate it with Find (including the trailing\n# space that must be removed) against a fresh Range built from the whole\n# document content, then overwrite just that matched range's Text with the\n# trimmed version, leaving the rest of the document untouched.\n\n$d = $word.ActiveDocument\n\n# 1. \"...dependiendo el tipo de trabajo realizado se tienen diferentes costos \"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.MatchCase = $true\n$found1 = $rng1.Find.Execute(\" dependiendo el tipo de trabajo realizado se tienen diferentes costos \")\nif ($found1) {\n    $rng1.Text = \" dependiendo el tipo de trabajo realizado se tienen diferentes costos\"\n}\n\n# 2. \"video conferencia \"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.MatchCase = $true\n$found2 = $rng2.Find.Execute(\"video conferencia \")\nif ($found2) {\n    $rng2.Text = \"video conferencia\"\n}\n\n# 3. \"8. \u00bfCu\u00e1l es su prop\u00f3sito? \" \u2014 drop the trailing space-only run after \"prop\u00f3sito?\"\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.MatchCase = $true\n$found3 = $rng3.Find.Execute(\"prop\u00f3sito? \")\nif ($found3) {\n    $rng3.Text = \"prop\u00f3sito?\"\n}\n\n# 4. \"Que necesito para realizar la consultor\u00eda? \"\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.MatchCase = $true\n$found4 = $rng4.Find.Execute(\"Que necesito para realizar la consultor\u00eda? \")\nif ($found4) {\n    $rng4.Text = \"Que necesito para realizar la consultor\u00eda?\"\n}\n\n# 5. \"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor? \"\n$rng5 = $d.Content\n$rng5.Find.ClearFormatting()\n$rng5.Find.MatchCase = $true\n$found5 = $rng5.Find.Execute(\"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor? \")\nif ($found5) {\n    $rng5.Text = \"\u00bfC\u00f3mo reservo las sesiones de trabajo con el consultor?\"\n}\n\n# 6. \"agenda otra cita \"\n$rng6 = $d.Content\n$rng6.Find.ClearFormatting()\n$rng6.Find.MatchCase = $true\n$found6 = $rng6.Find.Execute(\"agenda otra cita \")\nif ($found6) {\n    $rng6.Text = \"agenda otra cita\"\n}\n"}
